$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "alpha3F"

# Row 13 precision tweaks
$ws.Range("C13").Value = 0.9854415717608653
$ws.Range("D13").Value = 0.996290706389217
$ws.Range("F13").Value = 0.9854415717608653
$ws.Range("G13").Value = 0.9881389825467225
$ws.Range("H13").Value = 0.9998535884593858
$ws.Range("J13").Value = 0.996290706389217
$ws.Range("K13").Value = 0.9959943579930182

# Row 15 precision tweak
$ws.Range("H15").Value = 0.8684472120320456
